# Add new columns K:U to the contratos sheet, fill header row, create
# empty placeholder cells for existing data rows, and append a new
# data row (73) with values only in the new K:U columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells K1:U1 -------------------------------------------------
# Copy the formatting of the existing header cell (A1, style index 1: bold,
# bordered, centered) onto the new header range before writing values so the
# new headers look consistent with the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("K1:U1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K1").Value = "n_contrato"
$ws.Range("L1").Value = "n_licitacao"
$ws.Range("M1").Value = "assinatura"
$ws.Range("N1").Value = "vencimento"
$ws.Range("O1").Value = "contratada"
$ws.Range("P1").Value = "cnpj"
$ws.Range("Q1").Value = "modalidade"
$ws.Range("R1").Value = "objeto"
$ws.Range("S1").Value = "contratante"
$ws.Range("T1").Value = "valor"
$ws.Range("U1").Value = "filename"

# --- 2. Blank placeholder cells for K2:U72 -------------------------------------
# The existing rows (2-72) gain empty cells in the new columns so the sheet's
# rectangular extent covers A1:U72 before the new row is appended.
$ws.Range("K2:U72").Style = "Normal"

# --- 3. New row 73 --------------------------------------------------------------
# Columns A:J stay blank for this row; the new contract data lives in K:T,
# and U (filename) is left blank.
$ws.Range("A73:J73").Style = "Normal"
$ws.Range("U73").Style = "Normal"

$ws.Range("K73").Value = "154/2023"
$ws.Range("L73").Value = "068/2023"
$ws.Range("M73").Value = "19/12/2023"
$ws.Range("N73").Value = "31/12/2023"
$ws.Range("O73").Value = "LABORCOM COMÉRCIO DE MATERIAIS DE CONSTRUÇÃO LTDA."
$ws.Range("P73").Value = "34.101.659/0001-56"
$ws.Range("Q73").Value = "DISPENSA DE LICITAÇÃO"
$ws.Range("R73").Value = "Contratação de pessoa jurídica para fornecimento de material elétrico para ILUMINAÇÃO PÚBLICA em atendimento às necessidades da Secretaria de Infraestrutura e Urbanismo do Município de Nilo Peçanha - BA."
$ws.Range("S73").Value = "MUNICÍPIO DE NILO PEÇANHA"
$ws.Range("T73").Value = "54.720,00"
